# Apply "New sample reduction steps" update:
#  - Row 12/13 keep their data but get refreshed timestamps
#  - Row 14 becomes the new "binary_all_int_poly" step (new num_cols + timestamp)
#  - Old row 14 ("multi_all") and row 15 ("multi_nolags") content shift down to rows 15 and 16
#  - A brand new row 17 ("multi_all_int_poly") is appended
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - binary_all: only the timestamp changes
$ws.Cells.Item(12, 10).Value = 45016.44736511295

# Row 13 - binary_nolags: only the timestamp changes
$ws.Cells.Item(13, 10).Value = 45016.447382940634

# Row 14 - becomes binary_all_int_poly (new step inserted before old multi_all/multi_nolags)
$ws.Cells.Item(14, 1).Value = "estimation_sample"
$ws.Cells.Item(14, 2).Value = "binary_all_int_poly"
$ws.Cells.Item(14, 3).Value = "controls_same_outcome"
$ws.Cells.Item(14, 4).Value = "down"
$ws.Cells.Item(14, 5).Value = "weekly"
$ws.Cells.Item(14, 6).Value = "yes"
$ws.Cells.Item(14, 7).Value = 5218.0
$ws.Cells.Item(14, 8).Value = 11211.0
$ws.Cells.Item(14, 9).Value = 40804.0
$ws.Cells.Item(14, 10).Value = 45016.448366390265

# Row 15 - multi_all (shifted down from the old row 14)
$ws.Cells.Item(15, 1).Value = "estimation_sample"
$ws.Cells.Item(15, 2).Value = "multi_all"
$ws.Cells.Item(15, 3).Value = "controls_same_outcome"
$ws.Cells.Item(15, 4).Value = "down"
$ws.Cells.Item(15, 5).Value = "weekly"
$ws.Cells.Item(15, 6).Value = "yes"
$ws.Cells.Item(15, 7).Value = 5218.0
$ws.Cells.Item(15, 8).Value = 11211.0
$ws.Cells.Item(15, 9).Value = 537.0
$ws.Cells.Item(15, 10).Value = 45016.44838751201

# Row 16 - multi_nolags (shifted down from the old row 15) - new row
$ws.Cells.Item(16, 1).Value = "estimation_sample"
$ws.Cells.Item(16, 2).Value = "multi_nolags"
$ws.Cells.Item(16, 3).Value = "controls_same_outcome"
$ws.Cells.Item(16, 4).Value = "down"
$ws.Cells.Item(16, 5).Value = "weekly"
$ws.Cells.Item(16, 6).Value = "yes"
$ws.Cells.Item(16, 7).Value = 5218.0
$ws.Cells.Item(16, 8).Value = 11211.0
$ws.Cells.Item(16, 9).Value = 424.0
$ws.Cells.Item(16, 10).Value = 45016.44839460661

# Row 17 - brand new multi_all_int_poly step
$ws.Cells.Item(17, 1).Value = "estimation_sample"
$ws.Cells.Item(17, 2).Value = "multi_all_int_poly"
$ws.Cells.Item(17, 3).Value = "controls_same_outcome"
$ws.Cells.Item(17, 4).Value = "down"
$ws.Cells.Item(17, 5).Value = "weekly"
$ws.Cells.Item(17, 6).Value = "yes"
$ws.Cells.Item(17, 7).Value = 5218.0
$ws.Cells.Item(17, 8).Value = 11211.0
$ws.Cells.Item(17, 9).Value = 40808.0
$ws.Cells.Item(17, 10).Value = 45016.44968821267

# Apply the date/time number format to the new timestamp cells (same style as the rest of column J)
$ws.Range("J14:J17").NumberFormat = "m/d/yyyy h:mm:ss"
